$wb = $excel.ActiveWorkbook

# --- "Linear" sheet updates (mu, B, sig2, abs_epsi_autocorr) ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.0009214516252594898
$wsLinear.Range("B3").Value = -0.009188889892643243
$wsLinear.Range("B4").Value = 0.02079474206228022
$wsLinear.Range("B5").Value = "[1.0, 0.253769643185465, 0.09168133324142269, 0.06974536101673491, 0.06161771953519132, 0.10126183335825506, 0.2508525183222177, 0.3881695389017516, 0.2532700440636223, 0.09719352169554886, 0.06820010073429045, 0.07688238826077794, 0.09512319756070209, 0.2538640572007374, 0.4069881872758495, 0.24381236459766442, 0.056537286818766734, 0.06671164665678242, 0.07180881045748895, 0.08489772771777448]"

# --- "NonLinear" sheet updates (mu_0, B_0, sig2_0, mu_1, B_1, sig2_1, abs_epsi_autocorr) ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.008169737245117279
$wsNonLinear.Range("B5").Value = 0.1601177263490018
$wsNonLinear.Range("B6").Value = 0.01930805619532345
$wsNonLinear.Range("B7").Value = 0.0005750446363366753
$wsNonLinear.Range("B8").Value = -0.1011632280354307
$wsNonLinear.Range("B9").Value = 0.02226465185617424
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.2553927485275999, 0.10407769165931262, 0.08687993575737735, 0.07640305740193562, 0.11826643647463328, 0.25555697830038815, 0.3825389329050817, 0.254844714699065, 0.11020113727660379, 0.08332833375677075, 0.09228520452258848, 0.1083505429606133, 0.2558392713151685, 0.4015497169797526, 0.24675690976528838, 0.06939479329079168, 0.07942024365828708, 0.08317956420166256, 0.09618277635246343]"

$wb.Save()
